$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Source"
$ws.Range("B1").Value = "Amount"
$ws.Range("C1").Value = "Date"

# Row 2 - Salary income
$ws.Range("A2").Value = "Salary"
$ws.Range("B2").Value = 1000
$ws.Range("C2").Value = 45900.229537037034
$ws.Range("C2").NumberFormat = "mm-dd-yy"

# Row 3 - Interest income (trailing space in label is intentional)
$ws.Range("A3").Value = "Interest from saving account "
$ws.Range("B3").Value = 5000
$ws.Range("C3").Value = 45898.229537037034

# Re-use the exact same date style created for C2 (rather than creating a
# duplicate, functionally-identical style entry) by copying its format.
$ws.Range("C2").Copy()
$ws.Range("C3").PasteSpecial(-4122)
$excel.CutCopyMode = $false
